# The author re-uploaded the workbook after deleting one data row (row 84,
# the Al-Co-Cr-Fe-Ni composition entry with HV=733) from the sorted
# CITRINE hardness dataset sheet. Deleting the whole row shifts every
# subsequent row up by one, which is exactly what the diff shows (row 84's
# new contents equal old row 85's contents, row 85 <- old 86, ... row 103
# <- old 104, and old row 104 disappears). It also shrinks the used range
# from H104 to H103 and drops the shared-string cell count from 700 to 693
# (the deleted row contributed 7 text cells, columns A-G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 84 - Excel automatically shifts rows 85:104 up to
# fill the gap, shrinking the sheet's dimension to A1:H103.
$ws.Rows(84).Delete()

# Reflect the selection left behind in the file: the user had just deleted
# the row while it (the whole row) was selected, so the new row 84 (now
# holding what used to be row 85) ends up selected/active.
$ws.Range("A84:XFD84").Select()
